$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q1").Value = "Coupons"
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "PLU Lookup"
$ws.Range("Q20").Value = 8578

Write-Host "Q1:" $ws.Range("Q1").Value()
Write-Host "A20:" $ws.Range("A20").Value()
Write-Host "B20:" $ws.Range("B20").Value()
Write-Host "Q20:" $ws.Range("Q20").Value()
